# InputData.xlsx test-data refresh
#
# The QA automation suite re-generates the "used" login e-mail on each
# run (see juan.perez99_<random>@test.com in LoginData!A2) so the same
# address isn't reused across registration/login test passes.
# This commit just rolls that value forward to the latest generated one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")
$ws.Range("A2").Value = "juan.perez99_20113@test.com"
